$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each calendar-year block of 12 monthly rows is rotated so that the
# Oct/Nov/Dec rows move to the front of the block, followed by Jan..Sep.
# Year blocks start at row 2 (2014), 14 (2015), 26 (2016), 38 (2017).

$blockStarts = @(2, 14, 26, 38)

foreach ($start in $blockStarts) {
    $end = $start + 11

    # capture the current (pre-shift) values for this block of 12 rows, columns A:F
    $original = @{}
    for ($r = $start; $r -le $end; $r++) {
        $rowVals = @()
        for ($c = 1; $c -le 6; $c++) {
            $rowVals += , ($ws.Cells.Item($r, $c).Value())
        }
        $original[$r] = $rowVals
    }

    # write back rotated: new row r gets values that were 3 rows later (wrapping within block)
    for ($r = $start; $r -le $end; $r++) {
        $offset = $r - $start
        $newOffset = ($offset + 9) % 12
        $srcRow = $start + $newOffset
        $vals = $original[$srcRow]
        for ($c = 1; $c -le 6; $c++) {
            $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
        }
    }
}
